$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Prix Spot": insert a new date column ("10-nov") before the
# "01-oct." column (currently column DO), shifting all following columns
# (old DO:ES, i.e. the October data) one column to the right (new DP:ET).
# ----------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column at DO (shifts DO:ES -> DP:ET)
$wsSpot.Range("DO1:DO25").EntireColumn.Insert()

# New header cell for the inserted column
$wsSpot.Range("DO1").Value = "10-nov"

# The inserted column has no data yet for this date -> "-" placeholder,
# matching the other not-yet-available date columns on the sheet.
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 119).Value = "-"
}

# ----------------------------------------------------------------------
# Sheet "Gaz": append two new daily rows.
# ----------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A147").Value = "'2025-11-08"
$wsGaz.Range("A147").Style = "Normal"
$wsGaz.Range("B147").Value = 29.755

$wsGaz.Range("A148").Value = "'2025-11-09"
$wsGaz.Range("A148").Style = "Normal"
$wsGaz.Range("B148").Value = 29.755

# ----------------------------------------------------------------------
# Sheet "CO2": append two new daily rows.
# ----------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A147").Value = "'2025-11-08"
$wsCO2.Range("A147").Style = "Normal"
$wsCO2.Range("B147").Value = 79.36

$wsCO2.Range("A148").Value = "'2025-11-09"
$wsCO2.Range("A148").Style = "Normal"
$wsCO2.Range("B148").Value = 79.36
